$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "43.685.00"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "2.280.02"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "'114.09"
$ws.Range("E5").Value = "  +10.14%  "
$ws.Range("D6").Value = "'266.92"
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("E7").Value = "  +1.06%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").Value = "'48.11"
$ws.Range("E10").Value = "  +5.06%  "
$ws.Range("D11").Value = "'0.0935"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("E12").Value = "  +7.45%  "
$ws.Range("D13").Value = "'0.108"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "'15.61"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "2.621.56"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("D17").Value = "2.280.42"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "43.483.02"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("E19").Value = "  -1.42%  "
$ws.Range("D20").Value = "'7.03"
$ws.Range("E20").Value = "  +12.08%  "
$ws.Range("D21").Value = "'71.80"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("E22").Value = "  -4.50%  "
$ws.Range("D23").Value = "'9.87"
$ws.Range("E23").Value = "  +7.13%  "
$ws.Range("D24").Value = "'231.51"
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("E25").Value = "  -3.70%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'11.47"
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D28").Value = "'40.93"
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").Value = "'172.72"
$ws.Range("E31").Value = "  -2.83%  "
$ws.Range("D32").Value = "'21.43"
$ws.Range("E32").Value = "  -1.94%  "
$ws.Range("D33").Value = "'0.0911"
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("D34").Value = "'5.63"
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -6.13%  "
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("E38").Value = "  -6.17%  "
$ws.Range("D39").Value = "'3.75"
$ws.Range("E39").Value = "  +5.65%  "
$ws.Range("D40").Value = "'14.37"
$ws.Range("E40").Value = "  +17.80%  "
$ws.Range("D41").Value = "'74.20"
$ws.Range("E41").Value = "  +12.55%  "
$ws.Range("D42").Value = "'2.41"
$ws.Range("E42").Value = "  +3.63%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "'6.18"
$ws.Range("E44").Value = "  +15.50%  "
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").Value = "'1.37"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "'8.65"
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("D48").Value = "'0.0998"
$ws.Range("E48").Value = "  -1.78%  "
$ws.Range("E49").Value = "  +1.46%  "
$ws.Range("D50").Value = "'101.08"
$ws.Range("E50").Value = "  +1.86%  "
$ws.Range("D51").Value = "'0.453"
$ws.Range("E51").Value = "  +3.43%  "
